# TimeLog_ConnorPeper.xlsx - continue work log entry for row 9 (Week 3)
# and reflect the current selection/scroll position in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the hours logged for week 3 (row 9)
$ws.Range("E9").Value = 4

# Append additional activity notes to the week-3 log entry
$ws.Range("F9").Value = "Meeting with professor. Met with team to choose SCRUM master and Product Owner. Meeting with professor for Daily SCRUM (Took longer than expected). Created Frontend-Backend connection."

# Move the selection to where editing left off
$ws.Activate()
$ws.Range("H13").Select()
